$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before current row 5 ("video" group) to make room
# for the new "product" pid_ndcg-style row and push "video" rows down.
# Plan: current rows 2-4 (product) become the new "ad" group (rows 2-4),
# a new "product" group (pid_pass/pid_recall/pid_ndcg) occupies rows 5-7,
# and the original "video" group moves to rows 8-10.

# Step 1: Insert 3 blank rows at row 5 to push "video" rows (5-7) down to (8-10).
$ws.Range("A5:H7").EntireRow.Insert()

# Step 2: Update the "ad" group (originally "product", rows 2-4).
$ws.Range("A2").Value = "ad"
$ws.Range("B2").Value = "pid_pass"
$ws.Range("C2").Value = 0.021
$ws.Range("D2").Value = 0.056
$ws.Range("E2").Value = 0.079
$ws.Range("F2").Value = 0.109
$ws.Range("G2").Value = 0.156
$ws.Range("H2").Value = 0.214

$ws.Range("A3").Value = "ad"
$ws.Range("B3").Value = "pid_recall"
$ws.Range("C3").Value = 0.006613492063492063
$ws.Range("D3").Value = 0.01824642857142857
$ws.Range("E3").Value = 0.02749603174603174
$ws.Range("F3").Value = 0.03750039682539683
$ws.Range("G3").Value = 0.05507500000000005
$ws.Range("H3").Value = 0.07451111111111113

$ws.Range("A4").Value = "ad"
$ws.Range("B4").Value = "pid_ndcg"
$ws.Range("C4").Value = 0.021
$ws.Range("D4").Value = 0.02254369849090224
$ws.Range("E4").Value = 0.02471010461444407
$ws.Range("F4").Value = 0.02745757381377225
$ws.Range("G4").Value = 0.03377001744513262
$ws.Range("H4").Value = 0.04010448322283051

# Step 3: Fill the new "product" group (rows 5-7).
$ws.Range("A5").Value = "product"
$ws.Range("B5").Value = "pid_pass"
$ws.Range("C5").Value = 0.012
$ws.Range("D5").Value = 0.037
$ws.Range("E5").Value = 0.054
$ws.Range("F5").Value = 0.083
$ws.Range("G5").Value = 0.114
$ws.Range("H5").Value = 0.153

$ws.Range("A6").Value = "product"
$ws.Range("B6").Value = "pid_recall"
$ws.Range("C6").Value = 0.001819047619047619
$ws.Range("D6").Value = 0.006647619047619048
$ws.Range("E6").Value = 0.01078492063492063
$ws.Range("F6").Value = 0.01550476190476189
$ws.Range("G6").Value = 0.02309642857142856
$ws.Range("H6").Value = 0.03220357142857141

$ws.Range("A7").Value = "product"
$ws.Range("B7").Value = "pid_ndcg"
$ws.Range("C7").Value = 0.012
$ws.Range("D7").Value = 0.01302162251422344
$ws.Range("E7").Value = 0.01309264166566387
$ws.Range("F7").Value = 0.01345816976840837
$ws.Range("G7").Value = 0.01659217718367545
$ws.Range("H7").Value = 0.02002961693412586

# Step 4: The "video" group is now in rows 8-10 (shifted down by the insert);
# values are unchanged from before, so nothing further needed there.
